# Actualización automática hashcode vie ene 11 01:41:26 CET 2019
# Updates the hashcode values (column B) in the "hashcode.csv" sheet to their
# new values, identified by matching the current (old) hash text per cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$map = @{
    "B11"  = @("27727890b1e1b6a6913a83c5b04b29a4", "8d34570a89896d1e7487e15264d4430d")
    "B24"  = @("0148032060184c38150a70dc4cc62862", "474cf07313af8791624506f75e5f63ca")
    "B29"  = @("283a43421d3d619311cfd0592b2dd6eb", "c9ab33bf5bace551342158f57f5fb0c5")
    "B121" = @("a75a4d9af7223344b490c2aca5cdac25", "27ce3918723a74c22be7d3b4776af7d0")
    "B191" = @("3a13e7d435e81d1a9016877eee3af917", "abec732590cbb485771a4e1b3fca3502")
    "B227" = @("79d7ac27c02b8ee4b146a8ebaf9cdac1", "366679d9cd102f7c634ebffd2d642faa")
    "B232" = @("ae22bcdb5a3d16e8e1bb7667b80435a8", "2ad3ae0d1889ca9238638c3c5377ba7a")
    "B419" = @("afba4ee92bb44bede48ddf483ac24705", "2ee5add6736bc97726d8045230c25adb")
    "B478" = @("19b25a4ce25f6f97839a85d363ab88b0", "0e421a028fe726870a018a31b7132a98")
    "B480" = @("54047bec7956934d2f51b05c58bf2b32", "18e444a0140e6b442fd1939ef7a91505")
    "B502" = @("81629ac93065ab0b8af54b4a0aeeec72", "2be9b891f2e904e9681becfeefd8ad95")
    "B629" = @("8e135f17d024197e1fee484b3eb87bd1", "00d68d50c3de3d47c92bdab22d9dc903")
    "B649" = @("3e72d49f2d5a1c266973b510c1bc866a", "759613b2f4e599e5bbf90a4d43e40cc9")
    "B655" = @("d6d55401dea2dc036bcb028447293785", "6a5e3c6b8da31df5f747f3f32e2ebcf8")
    "B665" = @("1ba24c61578dfbe6dd75691af4a3de32", "8d7ab49717672ca1dce808920279e201")
    "B680" = @("dfc9b3ba408aa959d34138ce25d08e59", "c56276a6b66cd48ad0785d014eb73047")
    "B685" = @("5bc27490b7119c501eed5d24ed5b0700", "ec3c9dece34deb785b2e0c3199198bf1")
    "B703" = @("19cfb9675ed60fea946f78fdbeb00be0", "638699b23a2f1096340b7b6ea891a02e")
    "B704" = @("aabab779119ff412ea1a22091217f45a", "6ce535b8d351b4b4aa19d5896b319229")
    "B742" = @("751a358c7da9cc56db9b1cdcade1cc19", "fec6821cbe9c86068c0b2ce65f3d2782")
    "B874" = @("d878f735a89572d2273c1e98708e28dd", "c9c849f03081bb7a17b5eba5feebb7ea")
}

foreach ($addr in $map.Keys) {
    $old = $map[$addr][0]
    $new = $map[$addr][1]
    $cell = $ws.Range($addr)
    $current = $cell.Value2
    if ($current -eq $old) {
        $cell.Value = $new
    } else {
        Write-Host "WARNING: cell $addr did not match expected old value. Found: $current"
        $cell.Value = $new
    }
}
